# Apply the 2021-07-08 -> 2021-07-09 "as of" date refresh and refreshed
# Weight / Percent Change figures for the SQE_holdings model-holdings sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet ships protected (legacy password hash "D382"); temporarily
# unprotect so the cell writes below are allowed, then restore protection
# with the same settings (contents/objects/scenarios locked, row & column
# formatting left allowed) once all edits are in place.
$ws.Unprotect()

# --- Update the "as of" date in the confidentiality / disclosure banner ---
$bannerText = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + "`n" + "Model holdings provided as of 2021-07-09 for illustrative purposes only and are subject to change."
$ws.Range("A33").Value = $bannerText

# --- Refresh the Weight (D) and Percent Change (E) columns, rows 2-30 ---
$ws.Range("D2").Value = 0.01974717605567974
$ws.Range("E2").Value = 0.004024819721616391
$ws.Range("D3").Value = 0.01824249741731455
$ws.Range("E3").Value = -0.002392936710949867
$ws.Range("D4").Value = 0.07529053047448407
$ws.Range("E4").Value = 0.003794664278174098
$ws.Range("D5").Value = 0.05616819645840556
$ws.Range("E5").Value = -0.003234702163525283
$ws.Range("D6").Value = 0.07115341685935515
$ws.Range("E6").Value = 0.01305501256632224
$ws.Range("D7").Value = 0.01993684143026358
$ws.Range("E7").Value = 0.009581263307310017
$ws.Range("D8").Value = 0.03327362888115944
$ws.Range("E8").Value = 0.03249097472924167
$ws.Range("D9").Value = 0.02886165100095902
$ws.Range("E9").Value = 0.009012392039053907
$ws.Range("D10").Value = 0.02362944563972099
$ws.Range("E10").Value = 0.006096434509514026
$ws.Range("D11").Value = 0.02629574943337469
$ws.Range("E11").Value = -0.003434655675768394
$ws.Range("D12").Value = 0.02601501457337559
$ws.Range("E12").Value = 0.01380008679299882
$ws.Range("D13").Value = 0.04315429172852744
$ws.Range("E13").Value = 0.01114459722483829
$ws.Range("D14").Value = 0.02334570021822055
$ws.Range("E14").Value = 0.01083228019498095
$ws.Range("D15").Value = 0.04072205909160221
$ws.Range("E15").Value = 0.003962621244381292
$ws.Range("D16").Value = 0.02953691994570754
$ws.Range("E16").Value = 0.03199946998807479
$ws.Range("D17").Value = 0.04457919048710735
$ws.Range("E17").Value = 0.007428617736837895
$ws.Range("D18").Value = 0.1169265960378903
$ws.Range("E18").Value = 0.001874414245548239
$ws.Range("D19").Value = 0.02901518963752689
$ws.Range("E19").Value = 0.002303430243416615
$ws.Range("D20").Value = 0.02396738116824538
$ws.Range("E20").Value = 0.007411036163344153
$ws.Range("D21").Value = 0.02447902609539657
$ws.Range("E21").Value = 0.02535342913891814
$ws.Range("D22").Value = 0.01332399256451521
$ws.Range("E22").Value = 0.01748856126080311
$ws.Range("D23").Value = 0.0147705673659047
$ws.Range("E23").Value = 0.009171974522292903
$ws.Range("D24").Value = 0.03092900358392294
$ws.Range("E24").Value = 0.0003650167907722768
$ws.Range("D25").Value = 0.01109994025540701
$ws.Range("E25").Value = -0.003010577705451656
$ws.Range("D26").Value = 0.037025390624118
$ws.Range("E26").Value = 0.005268935236004468
$ws.Range("D27").Value = 0.02346792901517459
$ws.Range("E27").Value = 0.007004310344827624
$ws.Range("D28").Value = 0.05342467176224272
$ws.Range("E28").Value = 0.007861037149740069
$ws.Range("D29").Value = 0.04161800219439828
$ws.Range("E29").Value = 0.02453703703703702
$ws.Range("E30").Value = 0.008416476261196015

# Restore sheet protection (same password + options as the original file).
$ws.Protect("D382", $true, $true, $true, $false, $false, $true, $true)

